$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.791393
$ws.Range("H2").Value = 17.374179
$ws.Range("I2").Value = 0.2508946350719245
$ws.Range("J2").Value = 0.2508946350719244
$ws.Range("M2").Value = 5.440305333333334
$ws.Range("N2").Value = 16.320916
$ws.Range("O2").Value = 0.1071123097215924
$ws.Range("P2").Value = 0.1071123097215924
$ws.Range("Q2").Value = 31.50694622532934
$ws.Range("R2").Value = 283.562516027964
$ws.Range("S2").Value = 0.02687390385930988
$ws.Range("T2").Value = 0.02687390385930987
$ws.Range("G3").Value = 5.791393
$ws.Range("H3").Value = 17.374179
$ws.Range("I3").Value = 0.2508946350719245
$ws.Range("J3").Value = 0.2508946350719244
$ws.Range("O3").Value = 0.2078855279387566
$ws.Range("P3").Value = 0.2078855279387566
$ws.Range("Q3").Value = 61.149256951092
$ws.Range("R3").Value = 550.3433125598281
$ws.Range("S3").Value = 0.0521573636689287
$ws.Range("T3").Value = 0.0521573636689287
$ws.Range("G4").Value = 5.791393
$ws.Range("H4").Value = 17.374179
$ws.Range("I4").Value = 0.2508946350719245
$ws.Range("J4").Value = 0.2508946350719244
$ws.Range("M4").Value = 28.38327833333333
$ws.Range("N4").Value = 85.149835
$ws.Range("O4").Value = 0.5588286527093509
$ws.Range("P4").Value = 0.5588286527093509
$ws.Range("Q4").Value = 164.3787194567183
$ws.Range("R4").Value = 1479.408475110465
$ws.Range("S4").Value = 0.1402071108892478
$ws.Range("T4").Value = 0.1402071108892478
$ws.Range("G5").Value = 5.791393
$ws.Range("H5").Value = 17.374179
$ws.Range("I5").Value = 0.2508946350719245
$ws.Range("J5").Value = 0.2508946350719244
$ws.Range("M5").Value = 6.408436333333334
$ws.Range("N5").Value = 19.225309
$ws.Range("O5").Value = 0.1261735096303001
$ws.Range("P5").Value = 0.1261735096303001
$ws.Range("Q5").Value = 37.11377332181234
$ws.Range("R5").Value = 334.0239598963111
$ws.Range("S5").Value = 0.03165625665443808
$ws.Range("T5").Value = 0.03165625665443807
$ws.Range("I6").Value = 0.3683465988617928
$ws.Range("J6").Value = 0.3683465988617928
$ws.Range("M6").Value = 5.440305333333334
$ws.Range("N6").Value = 16.320916
$ws.Range("O6").Value = 0.1071123097215924
$ws.Range("P6").Value = 0.1071123097215924
$ws.Range("Q6").Value = 46.25637562674267
$ws.Range("R6").Value = 416.307380640684
$ws.Range("S6").Value = 0.03945445498217951
$ws.Range("T6").Value = 0.03945445498217951
$ws.Range("I7").Value = 0.3683465988617928
$ws.Range("J7").Value = 0.3683465988617928
$ws.Range("O7").Value = 0.2078855279387566
$ws.Range("P7").Value = 0.2078855279387566
$ws.Range("Q7").Value = 89.77521904525199
$ws.Range("R7").Value = 807.9769714072679
$ws.Range("S7").Value = 0.07657392716882921
$ws.Range("T7").Value = 0.07657392716882921
$ws.Range("I8").Value = 0.3683465988617928
$ws.Range("J8").Value = 0.3683465988617928
$ws.Range("M8").Value = 28.38327833333333
$ws.Range("N8").Value = 85.149835
$ws.Range("O8").Value = 0.5588286527093509
$ws.Range("P8").Value = 0.5588286527093509
$ws.Range("Q8").Value = 241.3297606773517
$ws.Range("R8").Value = 2171.967846096165
$ws.Range("S8").Value = 0.2058426335720074
$ws.Range("T8").Value = 0.2058426335720074
$ws.Range("I9").Value = 0.3683465988617928
$ws.Range("J9").Value = 0.3683465988617928
$ws.Range("M9").Value = 6.408436333333334
$ws.Range("N9").Value = 19.225309
$ws.Range("O9").Value = 0.1261735096303001
$ws.Range("P9").Value = 0.1261735096303001
$ws.Range("Q9").Value = 54.48794140256567
$ws.Range("R9").Value = 490.391472623091
$ws.Range("S9").Value = 0.04647558313877669
$ws.Range("T9").Value = 0.04647558313877669
$ws.Range("G10").Value = 2.834746
$ws.Range("H10").Value = 8.504238000000001
$ws.Range("I10").Value = 0.1228068209481894
$ws.Range("J10").Value = 0.1228068209481894
$ws.Range("M10").Value = 5.440305333333334
$ws.Range("N10").Value = 16.320916
$ws.Range("O10").Value = 0.1071123097215924
$ws.Range("P10").Value = 0.1071123097215924
$ws.Range("Q10").Value = 15.42188378244534
$ws.Range("R10").Value = 138.796954042008
$ws.Range("S10").Value = 0.01315412224132661
$ws.Range("T10").Value = 0.01315412224132661
$ws.Range("G11").Value = 2.834746
$ws.Range("H11").Value = 8.504238000000001
$ws.Range("I11").Value = 0.1228068209481894
$ws.Range("J11").Value = 0.1228068209481894
$ws.Range("O11").Value = 0.2078855279387566
$ws.Range("P11").Value = 0.2078855279387566
$ws.Range("Q11").Value = 29.931073844424
$ws.Range("R11").Value = 269.379664599816
$ws.Range("S11").Value = 0.02552976080729472
$ws.Range("T11").Value = 0.02552976080729472
$ws.Range("G12").Value = 2.834746
$ws.Range("H12").Value = 8.504238000000001
$ws.Range("I12").Value = 0.1228068209481894
$ws.Range("J12").Value = 0.1228068209481894
$ws.Range("M12").Value = 28.38327833333333
$ws.Range("N12").Value = 85.149835
$ws.Range("O12").Value = 0.5588286527093509
$ws.Range("P12").Value = 0.5588286527093509
$ws.Range("Q12").Value = 80.45938472230334
$ws.Range("R12").Value = 724.13446250073
$ws.Range("S12").Value = 0.0686279702939952
$ws.Range("T12").Value = 0.06862797029399519
$ws.Range("G13").Value = 2.834746
$ws.Range("H13").Value = 8.504238000000001
$ws.Range("I13").Value = 0.1228068209481894
$ws.Range("J13").Value = 0.1228068209481894
$ws.Range("M13").Value = 6.408436333333334
$ws.Range("N13").Value = 19.225309
$ws.Range("O13").Value = 0.1261735096303001
$ws.Range("P13").Value = 0.1261735096303001
$ws.Range("Q13").Value = 18.16628926217134
$ws.Range("R13").Value = 163.496603359542
$ws.Range("S13").Value = 0.01549496760557291
$ws.Range("T13").Value = 0.01549496760557291
$ws.Range("G14").Value = 5.954296666666667
$ws.Range("H14").Value = 17.86289
$ws.Range("I14").Value = 0.2579519451180933
$ws.Range("J14").Value = 0.2579519451180933
$ws.Range("M14").Value = 5.440305333333334
$ws.Range("N14").Value = 16.320916
$ws.Range("O14").Value = 0.1071123097215924
$ws.Range("P14").Value = 0.1071123097215924
$ws.Range("Q14").Value = 32.39319191191556
$ws.Range("R14").Value = 291.53872720724
$ws.Range("S14").Value = 0.02762982863877642
$ws.Range("T14").Value = 0.02762982863877641
$ws.Range("G15").Value = 5.954296666666667
$ws.Range("H15").Value = 17.86289
$ws.Range("I15").Value = 0.2579519451180933
$ws.Range("J15").Value = 0.2579519451180933
$ws.Range("O15").Value = 0.2078855279387566
$ws.Range("P15").Value = 0.2078855279387566
$ws.Range("Q15").Value = 62.86929877372
$ws.Range("R15").Value = 565.82368896348
$ws.Range("S15").Value = 0.05362447629370399
$ws.Range("T15").Value = 0.05362447629370399
$ws.Range("G16").Value = 5.954296666666667
$ws.Range("H16").Value = 17.86289
$ws.Range("I16").Value = 0.2579519451180933
$ws.Range("J16").Value = 0.2579519451180933
$ws.Range("M16").Value = 28.38327833333333
$ws.Range("N16").Value = 85.149835
$ws.Range("O16").Value = 0.5588286527093509
$ws.Range("P16").Value = 0.5588286527093509
$ws.Range("Q16").Value = 169.0024595692389
$ws.Range("R16").Value = 1521.02213612315
$ws.Range("S16").Value = 0.1441509379541005
$ws.Range("T16").Value = 0.1441509379541005
$ws.Range("G17").Value = 5.954296666666667
$ws.Range("H17").Value = 17.86289
$ws.Range("I17").Value = 0.2579519451180933
$ws.Range("J17").Value = 0.2579519451180933
$ws.Range("M17").Value = 6.408436333333334
$ws.Range("N17").Value = 19.225309
$ws.Range("O17").Value = 0.1261735096303001
$ws.Range("P17").Value = 0.1261735096303001
$ws.Range("Q17").Value = 38.15773109811223
$ws.Range("R17").Value = 343.41957988301
$ws.Range("S17").Value = 0.03254670223151237
$ws.Range("T17").Value = 0.03254670223151237
